$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update price on row 2 (Owain / 2021-04-25)
$ws.Range("C2").Value = 950

# 2. Remove the two extra "Owain" rows (old rows 3 and 4); remaining rows shift up.
$ws.Rows("3:4").Delete()

# 3. Stage the Name/Price/Country/Vat (B:E) blocks from their current (post-delete) rows
#    into scratch rows, keyed to the final row order, so the 10-way reshuffle doesn't clobber data.
$ws.Range("B9:E9").Copy($ws.Range("B100"))
$ws.Range("B6:E6").Copy($ws.Range("B101"))
$ws.Range("B7:E7").Copy($ws.Range("B102"))
$ws.Range("B3:E3").Copy($ws.Range("B103"))
$ws.Range("B4:E4").Copy($ws.Range("B104"))
$ws.Range("B5:E5").Copy($ws.Range("B105"))
$ws.Range("B10:E10").Copy($ws.Range("B106"))
$ws.Range("B11:E11").Copy($ws.Range("B107"))
$ws.Range("B12:E12").Copy($ws.Range("B108"))
$ws.Range("B8:E8").Copy($ws.Range("B109"))

# 4. Copy the staged blocks back into the final row order.
$ws.Range("B100:E100").Copy($ws.Range("B3"))
$ws.Range("B101:E101").Copy($ws.Range("B4"))
$ws.Range("B102:E102").Copy($ws.Range("B5"))
$ws.Range("B103:E103").Copy($ws.Range("B6"))
$ws.Range("B104:E104").Copy($ws.Range("B7"))
$ws.Range("B105:E105").Copy($ws.Range("B8"))
$ws.Range("B106:E106").Copy($ws.Range("B9"))
$ws.Range("B107:E107").Copy($ws.Range("B10"))
$ws.Range("B108:E108").Copy($ws.Range("B11"))
$ws.Range("B109:E109").Copy($ws.Range("B12"))

# 5. Fix up Price and Country per the target layout (Name/Vat Number came along with the copy above).
$ws.Range("C3").Value = 30000
$ws.Range("D3").Value = "United States of America"
$ws.Range("C4").Value = 2800
$ws.Range("D4").Value = "Switzerland"
$ws.Range("C5").Value = 4500
$ws.Range("D5").Value = "France"
$ws.Range("C6").Value = 950
$ws.Range("D6").Value = "United Kingdom"
$ws.Range("C7").Value = 99500
$ws.Range("D7").Value = "United States of America"
$ws.Range("C8").Value = 16000
$ws.Range("D8").Value = "United States of America"
$ws.Range("C9").Value = 24500
$ws.Range("D9").Value = "Switzerland"
$ws.Range("C10").Value = 3950
$ws.Range("D10").Value = "France"
$ws.Range("C11").Value = 2500
$ws.Range("D11").Value = "France"
$ws.Range("C12").Value = 15000
$ws.Range("D12").Value = "United Kingdom"

# 6. Clean up the scratch area.
$ws.Range("B100:E109").Clear()
